$p = $ppt.ActivePresentation

# --- Slide 1: set the (currently empty) title text ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Example Slides"

# --- Add 4 new slides, using the same slide master's custom layouts ---
$master = $p.SlideMaster

# Slide 2 -> "Title and Content" layout
$layout2 = $master.CustomLayouts.Item(2)
$s2 = $p.Slides.AddSlide(2, $layout2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Slide 1"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "TEST"

# Slide 3 -> "Name Card" layout
$layout12 = $master.CustomLayouts.Item(12)
$s3 = $p.Slides.AddSlide(3, $layout12)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Slide 2"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Test"

# Slide 4 -> "Vertical Title and Text" layout
$layout14 = $master.CustomLayouts.Item(14)
$s4 = $p.Slides.AddSlide(4, $layout14)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Slide 3"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Test"

# Slide 5 -> "Quote with Caption" layout
$layout11 = $master.CustomLayouts.Item(11)
$s5 = $p.Slides.AddSlide(5, $layout11)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Last Slide"
